# Economic Dashboard update - 2025-12-03
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Style changes: toggle the "yellow highlight" date style on a handful of
#    date cells. Style index 47 = plain date format, 48 = highlighted date
#    format. We copy formats from stable, untouched reference cells so the
#    saved workbook reuses the existing cellXf entries instead of creating
#    new ones.
#    C3  is a stable cell using style 47 (no highlight)
#    C7  is a stable cell using style 48 (highlighted)
# ---------------------------------------------------------------------------

# Cells that should gain the highlight (47 -> 48)
# NOTE: PasteSpecial only honors the first area of a multi-area (Union)
# range in this environment, so each target cell is pasted individually.
$highlightTargets = @("N5","C32","C33","C34","N41","N42","N43","N44")
foreach ($ref in $highlightTargets) {
    $ws.Range("C7").Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}

# Cells that should lose the highlight (48 -> 47)
$plainTargets = @("N13","N14","C28","C29","C30","C31","N51")
foreach ($ref in $plainTargets) {
    $ws.Range("C3").Copy() | Out-Null
    $ws.Range($ref).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Value changes
# ---------------------------------------------------------------------------

# Row 5 - ADP Total NonFarm Private
$ws.Range("N5").Value = 45962
$ws.Range("Q5").Value = -31000
$ws.Range("R5").Value = 47000
$ws.Range("S5").Value = -29000
$ws.Range("T5").Value = -3000
$ws.Range("U5").Value = 104000

# Row 29 - 5yr, 5yr Forward
$ws.Range("N29").Value = 45993
$ws.Range("Q29").Value = 2.18
$ws.Range("R29").Value = 2.17
$ws.Range("S29").Value = ""
$ws.Range("T29").Value = ""
$ws.Range("U29").Value = 2.17

# Row 30 - 10yr TIPS
$ws.Range("N30").Value = 45993
$ws.Range("Q30").Value = 2.24
$ws.Range("R30").Value = 2.24
$ws.Range("S30").Value = ""
$ws.Range("T30").Value = ""

# Row 32 - IP M/M
$ws.Range("C32").Value = 45901
$ws.Range("F32").Value = 0.0009661502023092794
$ws.Range("G32").Value = -0.002569125217784962
$ws.Range("H32").Value = 0.001607051307331187
$ws.Range("I32").Value = 0.004463379393190303
$ws.Range("J32").Value = -0.001489203276247131

# Row 33 - IP Y/Y
$ws.Range("C33").Value = 45901
$ws.Range("F33").Value = 0.01622608918688194
$ws.Range("G33").Value = 0.008952424004962659
$ws.Range("H33").Value = 0.01615692613305024
$ws.Range("I33").Value = 0.005292683410311731
$ws.Range("J33").Value = 0.001133220308735612

# Row 34 - Cap Util
$ws.Range("C34").Value = 45901
$ws.Range("F34").Value = 75.8665
$ws.Range("G34").Value = 75.8847
$ws.Range("H34").Value = 76.1724
$ws.Range("I34").Value = 76.1431
$ws.Range("J34").Value = 75.8982

# Row 41 - Export Prices M/M
$ws.Range("N41").Value = 45901
$ws.Range("Q41").Value = 0
$ws.Range("R41").Value = 0.0006544502617800152
$ws.Range("S41").Value = 0.003282994090610725
$ws.Range("T41").Value = 0.004617414248021312
$ws.Range("U41").Value = -0.006553079947575369

# Row 42 - Export Prices Y/Y
$ws.Range("N42").Value = 45901
$ws.Range("Q42").Value = 0.03801765105227423
$ws.Range("R42").Value = 0.03171390013495289
$ws.Range("S42").Value = 0.0241286863270779
$ws.Range("T42").Value = 0.02628032345013481
$ws.Range("U42").Value = 0.01881720430107515

# Row 43 - Import Prices M/M
$ws.Range("N43").Value = 45901
$ws.Range("Q43").Value = 0
$ws.Range("R43").Value = 0.0007082152974506872
$ws.Range("S43").Value = 0.00284090909090895
$ws.Range("T43").Value = -0.001418439716311948
$ws.Range("U43").Value = -0.004940014114325986

# Row 44 - Import Prices Y/Y
$ws.Range("N44").Value = 45901
$ws.Range("Q44").Value = 0.00283889283179564
$ws.Range("R44").Value = -0.000707213578500667
$ws.Range("S44").Value = -0.004231311706629215
$ws.Range("T44").Value = -0.005649717514124173
$ws.Range("U44").Value = -0.00353356890459364

# Row 47 - FFR
$ws.Range("N47").Value = 45992
$ws.Range("R47").Value = 3.89
$ws.Range("S47").Value = 3.89
$ws.Range("T47").Value = 3.89

# Row 48 - 2y UST
$ws.Range("N48").Value = 45992
$ws.Range("Q48").Value = 3.54
$ws.Range("S48").Value = ""
$ws.Range("T48").Value = 3.47
$ws.Range("U48").Value = ""

# Row 49 - 5y UST
$ws.Range("N49").Value = 45992
$ws.Range("Q49").Value = 3.67
$ws.Range("S49").Value = ""
$ws.Range("T49").Value = 3.59
$ws.Range("U49").Value = ""

# Row 50 - 10y UST
$ws.Range("N50").Value = 45992
$ws.Range("Q50").Value = 4.09
$ws.Range("S50").Value = ""
$ws.Range("T50").Value = 4.02
$ws.Range("U50").Value = ""

# Row 52 - BAA
$ws.Range("N52").Value = 45992
$ws.Range("Q52").Value = 5.87
$ws.Range("S52").Value = ""
$ws.Range("U52").Value = ""

Write-Host "Dashboard update applied."
